$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "PT100"
$ws.Range("B4").Value = 0.77083333333333337
$ws.Range("B4").NumberFormat = $ws.Range("B2").NumberFormat

$ws.Range("C4").Select()
